# Resume edit: "Personal Qualities and other Skills" bullet list update.
#
# The bullet that used to read "All above Junior level" is removed, and every
# following bullet shifts up one slot (HTML/CSS/Bootstrap+Material UI,
# JavaScript ES6, the React stack + React Toastify/Flip Move, jQuery,
# Firebase). The final bullet ("Firebase") paragraph loses its numbering /
# text and becomes a second blank spacer paragraph (matching the blank
# paragraph that already followed it), while the true trailing spacer
# paragraph is left untouched.
#
# We rebuild each affected paragraph's full contents (pPr + runs) via
# Range.InsertXML so the exact run/proofErr boundaries from the target are
# reproduced instead of letting same-format runs coalesce. Paragraph objects
# are captured up front (via .Next()) before any edits happen, since
# re-searching by text after editing would be ambiguous (several bullets
# share substrings once earlier bullets have already been rewritten).

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPrCambria = '<w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Calibri" w:hAnsi="Cambria"/></w:rPr>'

# Locate the anchor bullet, then walk forward to grab the five bullets after it.
$pJunior = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*All above Junior level*") {
        $pJunior = $para
        break
    }
}
if ($pJunior -eq $null) {
    throw "Could not find the 'All above Junior level' bullet paragraph"
}

$pHtml = $pJunior.Next()
$pJs = $pHtml.Next()
$pReact = $pJs.Next()
$pJquery = $pReact.Next()
$pFirebase = $pJquery.Next()

# --- "All above Junior level" -> "HTML, CSS, Bootstrap" + ", Material UI"
$xmlHtml = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="900" w:hanging="180"/>' +
    '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r><w:t>HTML, CSS, Bootstrap</w:t></w:r>' +
    '<w:r><w:t>, Material UI</w:t></w:r>' +
    '</w:p>'
$pJunior.Range.InsertXML($xmlHtml)

# --- "HTML, CSS, Bootstrap" / ", Material UI" -> "JavaScript ES6"
$xmlJs = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="900" w:hanging="180"/>' +
    '</w:pPr>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>JavaScript ES6</w:t></w:r>' +
    '</w:p>'
$pHtml.Range.InsertXML($xmlJs)

# --- "JavaScript ES6" -> React stack incl. "React Toastify" / "React Flip Move"
$xmlReact = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="900" w:hanging="180"/>' +
    '</w:pPr>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">  React JS, React Redux, React Routing, Styled Components, React Reveal, React </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPrCambria + '<w:t>Toastify</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrCambria + '<w:t>, React Flip Move</w:t></w:r>' +
    '</w:p>'
$pJs.Range.InsertXML($xmlReact)

# --- React stack -> "  jQuery"
$xmlJquery = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="900" w:hanging="180"/>' +
    '</w:pPr>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">  jQuery</w:t></w:r>' +
    '</w:p>'
$pReact.Range.InsertXML($xmlJquery)

# --- "jQuery" -> "  Firebase"
$xmlFirebase = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="900" w:hanging="180"/>' +
    '</w:pPr>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">  Firebase</w:t></w:r>' +
    '</w:p>'
$pJquery.Range.InsertXML($xmlFirebase)

# --- "Firebase" -> blank spacer paragraph (numbering removed, ind -> 1440)
$xmlBlank = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:spacing w:line="360" w:lineRule="auto"/>' +
        '<w:ind w:left="1440"/>' +
    '</w:pPr>' +
    '</w:p>'
$pFirebase.Range.InsertXML($xmlBlank)

Write-Output "done"
